$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the geographic (IP address) data in column D, rows 2-14.
# The header "IPAddress" in D1 is kept; only the data values are cleared.
$ws.Range("D2:D14").ClearContents()
$ws.Range("A1").Select() | Out-Null
